# Updated cryptos list on Tue May 16 09:20:08 UTC 2023 with GitHub Actions
# Refresh price/volume figures, and fix the swapped EthereumClassic /
# LidoDAOToken rows (25 & 26 in the ranking).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values are apostrophe-prefixed so Excel keeps them as
# plain text (matching the source data's inline-string cells) instead of
# silently re-typing them as numbers.

$ws.Range('D2').Value = '27.324.95'
$ws.Range('E2').Value = '  -1.45%  '
$ws.Range('D3').Value = '1.828.31'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  -0.79%  '
$ws.Range('D5').Value = '''314.34'
$ws.Range('E5').Value = '  -1.72%  '
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('D7').Value = '''0.4248'
$ws.Range('E7').Value = '  -1.84%  '
$ws.Range('D8').Value = '''0.3713'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('D9').Value = '''0.07262'
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('D10').Value = '''0.8666'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').Value = '''21.13'
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').Value = '1.844.53'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').Value = '''6.733'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').Value = '''0.07096'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').Value = '''5.324'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('D17').Value = '''1.007'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').Value = '''0.000008879'
$ws.Range('E18').Value = '  -1.78%  '
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').Value = '''15.10'
$ws.Range('E20').Value = '  -2.83%  '
$ws.Range('D21').Value = '27.327.20'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').Value = '''5.131'
$ws.Range('E22').Value = '  -2.56%  '
$ws.Range('D23').Value = '''10.91'
$ws.Range('E23').Value = '  -2.71%  '
$ws.Range('D24').Value = '2.049.37'
$ws.Range('E24').Value = '  -1.99%  '
$ws.Range('D25').Value = '''1.995'
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('D26').Value = '''153.01'
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''18.42'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''2.176'
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('D29').Value = '''5.245'
$ws.Range('E29').Value = '  -3.43%  '
$ws.Range('D30').Value = '''116.61'
$ws.Range('E30').Value = '  -3.37%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').Value = '''1.199'
$ws.Range('E32').Value = '  -3.28%  '
$ws.Range('D33').Value = '''0.7586'
$ws.Range('E33').Value = '  -2.79%  '
$ws.Range('D34').Value = '''4.465'
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('D35').Value = '''2.825'
$ws.Range('E35').Value = '  -3.25%  '
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('D37').Value = '''1.121'
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('D38').Value = '''0.01981'
$ws.Range('D39').Value = '''0.05271'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').Value = '''7.381'
$ws.Range('E40').Value = '  +3.02%  '
$ws.Range('D41').Value = '''2.869'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = '''0.1702'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('D43').Value = '''0.5067'
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('D44').Value = '''8.706'
$ws.Range('E44').Value = '  -2.89%  '
$ws.Range('D45').Value = '''10.62'
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('D46').Value = '''107.78'
$ws.Range('E46').Value = '  -2.83%  '
$ws.Range('D47').Value = '''0.4766'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('E49').Value = '  -2.62%  '
$ws.Range('D50').Value = '''0.06388'
$ws.Range('E50').Value = '  -1.89%  '
$ws.Range('D51').Value = '''1.860'
$ws.Range('E51').Value = '  -2.09%  '
